$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.123652100563049
$ws.Range("B1").Value = 3.424921274185181
$ws.Range("C1").Value = 4.526587009429932
$ws.Range("D1").Value = 2.086937665939331
$ws.Range("E1").Value = 1.576254963874817
